$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 59408
$ws.Range("C10").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D10").Value = 388.17
$ws.Range("E10").Value = 463.78
$ws.Range("F10").Value = 19
$ws.Range("G10").Value = 7375.23
$ws.Range("B11").Value = 47438
$ws.Range("C11").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 401.81
$ws.Range("E11").Value = 480.05
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 803.62
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 1330.35
$ws.Range("F26").Value = 114
$ws.Range("G26").Value = 5255.4
$ws.Range("B46").Value = 38469.98
$ws.Range("F76").Value = 16
$ws.Range("G76").Value = 1246.88
$ws.Range("B85").Value = 205652.14
$ws.Range("F120").Value = 321
$ws.Range("G120").Value = 26055.57
$ws.Range("B134").Value = 88499.42999999999
$ws.Range("F185").Value = 73
$ws.Range("G185").Value = 4766.9
$ws.Range("F191").Value = 75
$ws.Range("G191").Value = 9489.75
$ws.Range("F197").Value = 32
$ws.Range("G197").Value = 3345.6
$ws.Range("B198").Value = 60904.31
$ws.Range("F218").Value = 46
$ws.Range("G218").Value = 3599.96
$ws.Range("F225").Value = 59
$ws.Range("G225").Value = 4617.34
$ws.Range("B228").Value = 44627.96
$ws.Range("F248").Value = 82
$ws.Range("G248").Value = 1949.96
$ws.Range("F249").Value = 66
$ws.Range("G249").Value = 1901.46
$ws.Range("B267").Value = 35657.28
$ws.Range("F282").Value = 1
$ws.Range("G282").Value = 597.0700000000001
$ws.Range("F288").Value = 170
$ws.Range("G288").Value = 4209.2
$ws.Range("B298").Value = 136515.53
$ws.Range("F354").Value = 71
$ws.Range("G354").Value = 5681.42
$ws.Range("F357").Value = 25
$ws.Range("G357").Value = 993
$ws.Range("F359").Value = 19
$ws.Range("G359").Value = 3057.1
$ws.Range("F361").Value = 98
$ws.Range("G361").Value = 7223.58
$ws.Range("F367").Value = 31
$ws.Range("G367").Value = 4106.88
$ws.Range("F381").Value = 228
$ws.Range("G381").Value = 5271.36
$ws.Range("F392").Value = 31
$ws.Range("G392").Value = 2657.32
$ws.Range("F399").Value = 357
$ws.Range("G399").Value = 20955.9
$ws.Range("F400").Value = 2
$ws.Range("G400").Value = 433.26
$ws.Range("F402").Value = 80
$ws.Range("G402").Value = 4358.4
$ws.Range("F409").Value = 346
$ws.Range("G409").Value = 59280.18
$ws.Range("F410").Value = 53
$ws.Range("G410").Value = 8012.01
$ws.Range("F413").Value = 13
$ws.Range("G413").Value = 801.3200000000001
$ws.Range("F415").Value = 11
$ws.Range("G415").Value = 1454.53
$ws.Range("F418").Value = 147
$ws.Range("G418").Value = 8742.09
$ws.Range("F420").Value = 121
$ws.Range("G420").Value = 11208.23
$ws.Range("F421").Value = 399
$ws.Range("G421").Value = 15860.25
$ws.Range("F422").Value = 49
$ws.Range("G422").Value = 7039.34
$ws.Range("B423").Value = 292471.89
$ws.Range("F425").Value = 47
$ws.Range("G425").Value = 8628.26
$ws.Range("F426").Value = 69
$ws.Range("G426").Value = 12667.02
$ws.Range("F427").Value = 0
$ws.Range("G427").Value = 0
$ws.Range("B437").Value = 48940.58
$ws.Range("F464").Value = 26
$ws.Range("G464").Value = 6799.26
$ws.Range("B481").Value = 64847.22
$ws.Range("F488").Value = 47
$ws.Range("G488").Value = 500.55
$ws.Range("F491").Value = 38
$ws.Range("G491").Value = 1790.94
$ws.Range("F495").Value = 72
$ws.Range("G495").Value = 11951.28
$ws.Range("B497").Value = 60678.4
$ws.Range("F572").Value = 0
$ws.Range("G572").Value = 0
$ws.Range("B583").Value = 96436.87
$ws.Range("F595").Value = 117
$ws.Range("G595").Value = 5859.36
$ws.Range("F597").Value = 260
$ws.Range("G597").Value = 16029
$ws.Range("F600").Value = 85
$ws.Range("G600").Value = 2742.95
$ws.Range("F603").Value = 24
$ws.Range("G603").Value = 3823.68
$ws.Range("B610").Value = 86373
$ws.Range("F619").Value = 367
$ws.Range("G619").Value = 15736.96
$ws.Range("F621").Value = 312
$ws.Range("G621").Value = 18944.64
$ws.Range("B638").Value = 180788.88
$ws.Range("F666").Value = 58
$ws.Range("G666").Value = 4017.66
$ws.Range("B667").Value = 34668.98
$ws.Range("F674").Value = 124
$ws.Range("G674").Value = 6470.32
$ws.Range("F686").Value = 99
$ws.Range("G686").Value = 5393.52
$ws.Range("B688").Value = 108058.38
$ws.Range("F750").Value = 28
$ws.Range("G750").Value = 4247.6
$ws.Range("F767").Value = 8
$ws.Range("G767").Value = 3347.68
$ws.Range("B773").Value = 222781.75
$ws.Range("F800").Value = 23
$ws.Range("G800").Value = 1429.91
$ws.Range("F802").Value = 38
$ws.Range("G802").Value = 2136.36
$ws.Range("B805").Value = 8313.66
$ws.Range("F811").Value = 121
$ws.Range("G811").Value = 10611.7
$ws.Range("F814").Value = 117
$ws.Range("G814").Value = 16418.61
$ws.Range("B815").Value = 47420.95
$ws.Range("F817").Value = 79
$ws.Range("G817").Value = 6443.24
$ws.Range("F821").Value = 200
$ws.Range("G821").Value = 26620
$ws.Range("F824").Value = 79
$ws.Range("G824").Value = 17055.31
$ws.Range("F825").Value = 65
$ws.Range("G825").Value = 2423.2
$ws.Range("F830").Value = 358
$ws.Range("G830").Value = 54104.54
$ws.Range("F831").Value = 57
$ws.Range("G831").Value = 15587.22
$ws.Range("B837").Value = 248926.38
$ws.Range("F839").Value = 66
$ws.Range("G839").Value = 16465.68
$ws.Range("F840").Value = 55
$ws.Range("G840").Value = 9502.35
$ws.Range("F842").Value = 50
$ws.Range("G842").Value = 8194.5
$ws.Range("F847").Value = 58
$ws.Range("G847").Value = 2907.54
$ws.Range("F849").Value = 76
$ws.Range("G849").Value = 11121.08
$ws.Range("F852").Value = 47
$ws.Range("G852").Value = 6710.19
$ws.Range("F866").Value = 62
$ws.Range("G866").Value = 3533.38
$ws.Range("B867").Value = 272788.86
$ws.Range("F907").Value = 255
$ws.Range("G907").Value = 41593.05
$ws.Range("B911").Value = 42558.38
$ws.Range("B923").Value = 3913093.5
$ws.Range("B924").Value = 3913093.5

Write-Output "Applied changes"